$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    # Writing the literal string "2014-05-19" straight into .Value gets
    # "smart"-parsed as a date (like typing it into Excel would). Route it
    # through a text formula + copy/paste-values so the stored cell stays a
    # literal string, matching how the source data was corrected.
    $cell.Formula = "=""2014-05-19"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
